$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (Mandataire record) with new values
$ws.Range("A2").Value = "Ahmed Test"
$ws.Range("B2").Value = "BG1949"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "117165465787878754625432"
$ws.Range("D2").Value = "bmce"
$ws.Range("E2").Value = "bmce"
$ws.Range("G2").Value = "794/DR KESH"
$ws.Range("I2").Value = 10000
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 8500

# Update row 4 totals to reflect new sums
$ws.Range("I4").Value = 29000
$ws.Range("J4").Value = 2200
$ws.Range("K4").Value = 26800
